$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 418, shifting rows 418:521 down to 419:522
$ws.Rows.Item(418).Insert()

# Populate the newly inserted row 418 with the new data record
$ws.Cells.Item(418, 1).Value = 7
$ws.Cells.Item(418, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(418, 3).Value = "Ñuble"
$ws.Cells.Item(418, 4).Value = 44798
$ws.Cells.Item(418, 5).Value = 16
$ws.Cells.Item(418, 6).Value = 100112020
$ws.Cells.Item(418, 7).Value = "Tomate"
$ws.Cells.Item(418, 8).Value = "Larga vida"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 240
$ws.Cells.Item(418, 11).Value = 7500
$ws.Cells.Item(418, 12).Value = 8000
$ws.Cells.Item(418, 13).Value = 7750
$ws.Cells.Item(418, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(418, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(418, 16).Value = 431
$ws.Cells.Item(418, 17).Value = 18
$ws.Cells.Item(418, 18).Value = "Hortaliza"
